# Applies the updates described by the diff:
#  - Shared string "Random Forest" -> "Decision Tree" (cell B1, merged B1:Y1)
#  - Numeric values in columns L..X for rows 4..10 updated to new schedule counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the title text in B1 (merged B1:Y1) from "Random Forest" to "Decision Tree"
$ws.Range("B1").Value = "Decision Tree"

# 2) Update the numeric grid values (columns L..X) for rows 4 through 10
$data = @{
    4  = @{ L = 2; M = 2; N = 4; O = 5; P = 6; Q = 6; R = 6; S = 6; T = 6; U = 5; V = 5; W = 4; X = 0 }
    5  = @{ L = 1; M = 1; N = 3; O = 4; P = 5; Q = 5; R = 5; S = 5; T = 5; U = 5; V = 4; W = 3; X = 0 }
    6  = @{ L = 3; M = 3; N = 4; O = 6; P = 7; Q = 7; R = 7; S = 7; T = 7; U = 6; V = 5; W = 3; X = 0 }
    7  = @{ L = 3; M = 3; N = 4; O = 5; P = 7; Q = 7; R = 7; S = 7; T = 7; U = 5; V = 4; W = 2; X = 0 }
    8  = @{ L = 2; M = 2; N = 3; O = 5; P = 6; Q = 6; R = 6; S = 6; T = 5; U = 5; V = 5; W = 3; X = 0 }
    9  = @{ L = 1; M = 1; N = 3; O = 4; P = 6; Q = 6; R = 6; S = 6; T = 5; U = 5; V = 5; W = 3; X = 0 }
    10 = @{ L = 2; M = 2; N = 3; O = 4; P = 6; Q = 6; R = 6; S = 6; T = 5; U = 5; V = 5; W = 3; X = 0 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
